# Apply updates to the "想去人数" (F column) counts across the workbook's
# sheets, matching the regenerated data snapshot described by the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (cell, expectedOldValue, newValue)
$updates = @{
    "展览" = @(
        @{ Cell = "F5";  Old = 8183;  New = 8184 }
        @{ Cell = "F8";  Old = 2147;  New = 2148 }
        @{ Cell = "F12"; Old = 1176;  New = 1177 }
        @{ Cell = "F15"; Old = 752;   New = 753 }
        @{ Cell = "F21"; Old = 6963;  New = 6964 }
        @{ Cell = "F23"; Old = 54433; New = 54440 }
        @{ Cell = "F24"; Old = 4282;  New = 4283 }
        @{ Cell = "F26"; Old = 1026;  New = 1027 }
        @{ Cell = "F30"; Old = 864;   New = 865 }
        @{ Cell = "F33"; Old = 2057;  New = 2468 }
        @{ Cell = "F37"; Old = 849;   New = 850 }
        @{ Cell = "F38"; Old = 1128;  New = 1130 }
        @{ Cell = "F41"; Old = 177;   New = 178 }
        @{ Cell = "F44"; Old = 133;   New = 134 }
        @{ Cell = "F46"; Old = 130;   New = 131 }
    )
    "演出" = @(
        @{ Cell = "F12"; Old = 96; New = 97 }
    )
    "本地生活" = @(
        @{ Cell = "F10"; Old = 1612; New = 1613 }
    )
    "全部类型" = @(
        @{ Cell = "F4";  Old = 8183;  New = 8184 }
        @{ Cell = "F6";  Old = 1612;  New = 1613 }
        @{ Cell = "F9";  Old = 2147;  New = 2148 }
        @{ Cell = "F14"; Old = 1176;  New = 1177 }
        @{ Cell = "F18"; Old = 6963;  New = 6964 }
        @{ Cell = "F20"; Old = 54434; New = 54440 }
        @{ Cell = "F25"; Old = 4282;  New = 4283 }
        @{ Cell = "F32"; Old = 2057;  New = 2468 }
        @{ Cell = "F34"; Old = 849;   New = 850 }
        @{ Cell = "F35"; Old = 1128;  New = 1130 }
        @{ Cell = "F42"; Old = 133;   New = 134 }
        @{ Cell = "F44"; Old = 130;   New = 131 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($change in $updates[$sheetName]) {
        $ws.Range($change.Cell).Value = $change.New
    }
}
